# Fixed names in manual_typing
#
# The "manual_type" column (C) used a few inconsistent/abbreviated category
# labels. Rename them to clearer names, matching whole cell contents only
# so nothing else on the sheet is touched:
#   "Other AA"     -> "Other AAs"
#   "AA-ish"       -> "AA-like"
#   "NB or friend" -> "Nucleobases +"
# ("MAA" is left as-is.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("manual_typing")

$col = $ws.Range("C1:C152")

# xlWhole = 2 (match entire cell contents, not a partial/substring match)
$col.Replace("Other AA", "Other AAs", 2)
$col.Replace("AA-ish", "AA-like", 2)
$col.Replace("NB or friend", "Nucleobases +", 2)

# Update the sheet's last-saved selection.
$ws.Range("G42").Select()
